$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 209, shifting existing rows 209-334 down to 210-335
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with the new data point
$ws.Cells.Item(209, 1).Value = 8
$ws.Cells.Item(209, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(209, 3).Value = "Coquimbo"
$ws.Cells.Item(209, 4).Value = 44879
$ws.Cells.Item(209, 5).Value = 4
$ws.Cells.Item(209, 6).Value = 100112012
$ws.Cells.Item(209, 7).Value = "Espinaca"
$ws.Cells.Item(209, 8).Value = "Sin especificar"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 1100
$ws.Cells.Item(209, 11).Value = 500
$ws.Cells.Item(209, 12).Value = 600
$ws.Cells.Item(209, 13).Value = 550
$ws.Cells.Item(209, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(209, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(209, 16).Value = 1100
$ws.Cells.Item(209, 17).Value = 0.5
$ws.Cells.Item(209, 18).Value = "Hortaliza"

# Apply the date-format number format (same as the other D-column cells) to the new row's D cell
$ws.Cells.Item(209, 4).NumberFormat = $ws.Cells.Item(210, 4).NumberFormat
